# ULOHS deck: add physical implementation rows for XKSTIZ(1..5) and
# XKSTOZ(1..5) to the ARC reservoir "variables" worksheet (rows 94-103),
# mirroring the existing block of rows (3-93) for the D/E/F/G/H/I/J and
# L/M/N and Q/R/S column groups.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each new row's "G" (and mirrored "L"/"Q") label, and whether the
# "Implemented?" flag (column H, mirrored via the I/J note columns) is a
# hard numeric value (26, formatted as scientific notation like the rest
# of the sheet) or "N" (not implemented, with an explanatory note in J).
$newRows = @(
    @{ Row = 94;  Name = "XKSTIZ(1)"; Implemented = $true  },
    @{ Row = 95;  Name = "XKSTIZ(2)"; Implemented = $true  },
    @{ Row = 96;  Name = "XKSTIZ(3)"; Implemented = $false },
    @{ Row = 97;  Name = "XKSTIZ(4)"; Implemented = $false },
    @{ Row = 98;  Name = "XKSTIZ(5)"; Implemented = $true  },
    @{ Row = 99;  Name = "XKSTOZ(1)"; Implemented = $true  },
    @{ Row = 100; Name = "XKSTOZ(2)"; Implemented = $true  },
    @{ Row = 101; Name = "XKSTOZ(3)"; Implemented = $false },
    @{ Row = 102; Name = "XKSTOZ(4)"; Implemented = $false },
    @{ Row = 103; Name = "XKSTOZ(5)"; Implemented = $true  }
)

foreach ($nr in $newRows) {
    $r = $nr.Row
    $name = $nr.Name

    # ARC design columns (D, E, F) - same for every row in this block.
    $ws.Range("D$r").Value = "Y"
    $ws.Range("E$r").Value = "N"
    $ws.Range("F$r").Value = "N"

    # Variable name + old/new value columns (G, H, I, J).
    $ws.Range("G$r").Value = $name

    if ($nr.Implemented) {
        $ws.Range("H$r").Value = 26
        $ws.Range("H$r").NumberFormat = "0.00E+00"
    } else {
        $ws.Range("H$r").Value = "-"
    }

    $ws.Range("I$r").Value = 26
    $ws.Range("I$r").NumberFormat = "0.00E+00"
    $ws.Range("I$r").Font.Color = 393372

    if (-not $nr.Implemented) {
        $ws.Range("J$r").Value = "assumed the same as other zones"
    }

    # Mirrored "old" block (L, M, N).
    $ws.Range("L$r").Value = $name
    $ws.Range("M$r").Value = "-"
    $ws.Range("N$r").Value = "-"

    # Mirrored "new" block (Q, R, S).
    $ws.Range("Q$r").Value = $name
    $ws.Range("R$r").Value = "-"
    $ws.Range("S$r").Value = "-"
}

# Extend the "contains N" conditional-formatting highlight rules on
# columns D, E, F so they keep covering the whole table through row 103.
$cfD = $ws.Range("D3").FormatConditions.Item(1)
$cfD.ModifyAppliesToRange($ws.Range("D3:D103"))

$cfE = $ws.Range("E3").FormatConditions.Item(1)
$cfE.ModifyAppliesToRange($ws.Range("E3:E103"))

$cfF = $ws.Range("F3").FormatConditions.Item(1)
$cfF.ModifyAppliesToRange($ws.Range("F3:F103"))
